# Applies the 15 new training/match rows (1303-1317) plus related
# shared-string / dimension / selection updates described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy number/alignment formats for the styled columns (B = date style,
# D = centered shared-string style) from the last pre-existing data row so
# the new rows reuse the workbook's existing style indices instead of
# minting new ones. ---
$ws.Range("B1297").Copy()
$ws.Range("B1303:B1317").PasteSpecial(-4122)
$ws.Range("D1297").Copy()
$ws.Range("D1303:D1317").PasteSpecial(-4122)
$ws.Range("F245").Copy()
$ws.Range("F1317").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 1303
$ws.Range("A1303").Value = "Entrainement"
$ws.Range("B1303").Value = 46066
$ws.Range("C1303").Value = "Global"
$ws.Range("D1303").Value = "J-1"
$ws.Range("E1303").Value = "Mattheo Haon"
$ws.Range("F1303").Value = "right back"
$ws.Range("G1303").Value = "01:27:06"
$ws.Range("H1303").Value = 6.26
$ws.Range("I1303").Value = 0.51
$ws.Range("J1303").Value = 5.74
$ws.Range("K1303").Value = 0.36
$ws.Range("L1303").Value = 0.14
$ws.Range("M1303").Value = 0.03
$ws.Range("N1303").Value = 0
$ws.Range("O1303").Value = 3
$ws.Range("P1303").Value = 4.21
$ws.Range("Q1303").Value = 26.91
$ws.Range("R1303").Value = 4.72
$ws.Range("S1303").Value = 37
$ws.Range("T1303").Value = 9
$ws.Range("U1303").Value = 18
$ws.Range("V1303").Value = 9

# Row 1304
$ws.Range("A1304").Value = "Entrainement"
$ws.Range("B1304").Value = 46066
$ws.Range("C1304").Value = "Global"
$ws.Range("D1304").Value = "J-1"
$ws.Range("E1304").Value = "Karahali Souaré"
$ws.Range("F1304").Value = "right forward"
$ws.Range("G1304").Value = "01:26:47"
$ws.Range("H1304").Value = 5.18
$ws.Range("I1304").Value = 0.13
$ws.Range("J1304").Value = 5.04
$ws.Range("K1304").Value = 0.12
$ws.Range("L1304").Value = 0.02
$ws.Range("M1304").Value = 0
$ws.Range("N1304").Value = 0
$ws.Range("O1304").Value = 0
$ws.Range("P1304").Value = 3.3
$ws.Range("Q1304").Value = 24.31
$ws.Range("R1304").Value = 5.8
$ws.Range("S1304").Value = 42
$ws.Range("T1304").Value = 18
$ws.Range("U1304").Value = 33
$ws.Range("V1304").Value = 12

# Row 1305
$ws.Range("A1305").Value = "Entrainement"
$ws.Range("B1305").Value = 46066
$ws.Range("C1305").Value = "Global"
$ws.Range("D1305").Value = "J-1"
$ws.Range("E1305").Value = "Ilan Ihaddadene"
$ws.Range("F1305").Value = "center midfield"
$ws.Range("G1305").Value = "01:26:53"
$ws.Range("H1305").Value = 5.81
$ws.Range("I1305").Value = 0.21
$ws.Range("J1305").Value = 5.59
$ws.Range("K1305").Value = 0.18
$ws.Range("L1305").Value = 0.03
$ws.Range("M1305").Value = 0
$ws.Range("N1305").Value = 0
$ws.Range("O1305").Value = 0
$ws.Range("P1305").Value = 3.94
$ws.Range("Q1305").Value = 21.95
$ws.Range("R1305").Value = 5.45
$ws.Range("S1305").Value = 26
$ws.Range("T1305").Value = 11
$ws.Range("U1305").Value = 6
$ws.Range("V1305").Value = 2

# Row 1306
$ws.Range("A1306").Value = "N3 J16 Bourgoin 14/02/2026"
$ws.Range("B1306").Value = 46067
$ws.Range("C1306").Value = "Global"
$ws.Range("D1306").Value = "M"
$ws.Range("E1306").Value = "Mattheo Haon"
$ws.Range("F1306").Value = "right back"
$ws.Range("G1306").Value = "01:38:02"
$ws.Range("H1306").Value = 11.31
$ws.Range("I1306").Value = 2.09
$ws.Range("J1306").Value = 9.19
$ws.Range("K1306").Value = 1.3
$ws.Range("L1306").Value = 0.56
$ws.Range("M1306").Value = 0.24
$ws.Range("N1306").Value = 0.01
$ws.Range("O1306").Value = 18
$ws.Range("P1306").Value = 6.87
$ws.Range("Q1306").Value = 30.78
$ws.Range("R1306").Value = 4.74
$ws.Range("S1306").Value = 41
$ws.Range("T1306").Value = 13
$ws.Range("U1306").Value = 41
$ws.Range("V1306").Value = 12

# Row 1307
$ws.Range("A1307").Value = "N3 J16 Bourgoin 14/02/2027"
$ws.Range("B1307").Value = 46067
$ws.Range("C1307").Value = "Global"
$ws.Range("D1307").Value = "M"
$ws.Range("E1307").Value = "Yoan Zouma"
$ws.Range("F1307").Value = "center back"
$ws.Range("G1307").Value = "01:39:28"
$ws.Range("H1307").Value = 9.7
$ws.Range("I1307").Value = 1.14
$ws.Range("J1307").Value = 8.54
$ws.Range("K1307").Value = 0.81
$ws.Range("L1307").Value = 0.29
$ws.Range("M1307").Value = 0.06
$ws.Range("N1307").Value = 0
$ws.Range("O1307").Value = 8
$ws.Range("P1307").Value = 5.72
$ws.Range("Q1307").Value = 28.7
$ws.Range("R1307").Value = 4.49
$ws.Range("S1307").Value = 27
$ws.Range("T1307").Value = 6
$ws.Range("U1307").Value = 21
$ws.Range("V1307").Value = 12

# Row 1308
$ws.Range("A1308").Value = "N3 J16 Bourgoin 14/02/2028"
$ws.Range("B1308").Value = 46067
$ws.Range("C1308").Value = "Global"
$ws.Range("D1308").Value = "M"
$ws.Range("E1308").Value = "Ilan Ihaddadene"
$ws.Range("F1308").Value = "center midfield"
$ws.Range("G1308").Value = "00:24:28"
$ws.Range("H1308").Value = 3.06
$ws.Range("I1308").Value = 0.82
$ws.Range("J1308").Value = 2.23
$ws.Range("K1308").Value = 0.53
$ws.Range("L1308").Value = 0.24
$ws.Range("M1308").Value = 0.06
$ws.Range("N1308").Value = 0
$ws.Range("O1308").Value = 5
$ws.Range("P1308").Value = 7.46
$ws.Range("Q1308").Value = 28
$ws.Range("R1308").Value = 4.77
$ws.Range("S1308").Value = 10
$ws.Range("T1308").Value = 3
$ws.Range("U1308").Value = 10
$ws.Range("V1308").Value = 2

# Row 1309
$ws.Range("A1309").Value = "N3 J16 Bourgoin 14/02/2029"
$ws.Range("B1309").Value = 46067
$ws.Range("C1309").Value = "Global"
$ws.Range("D1309").Value = "M"
$ws.Range("E1309").Value = "Kamal Bafounta"
$ws.Range("F1309").Value = "center midfield"
$ws.Range("G1309").Value = "01:39:37"
$ws.Range("H1309").Value = 12.63
$ws.Range("I1309").Value = 3.26
$ws.Range("J1309").Value = 9.34
$ws.Range("K1309").Value = 2.32
$ws.Range("L1309").Value = 0.77
$ws.Range("M1309").Value = 0.19
$ws.Range("N1309").Value = 0.02
$ws.Range("O1309").Value = 13
$ws.Range("P1309").Value = 7.53
$ws.Range("Q1309").Value = 30.88
$ws.Range("R1309").Value = 4.65
$ws.Range("S1309").Value = 55
$ws.Range("T1309").Value = 6
$ws.Range("U1309").Value = 60
$ws.Range("V1309").Value = 9

# Row 1310
$ws.Range("A1310").Value = "N3 J16 Bourgoin 14/02/2030"
$ws.Range("B1310").Value = 46067
$ws.Range("C1310").Value = "Global"
$ws.Range("D1310").Value = "M"
$ws.Range("E1310").Value = "Sofiane Belle"
$ws.Range("F1310").Value = "left forward"
$ws.Range("G1310").Value = "01:24:46"
$ws.Range("H1310").Value = 9.26
$ws.Range("I1310").Value = 1.81
$ws.Range("J1310").Value = 7.42
$ws.Range("K1310").Value = 1.06
$ws.Range("L1310").Value = 0.63
$ws.Range("M1310").Value = 0.14
$ws.Range("N1310").Value = 0
$ws.Range("O1310").Value = 14
$ws.Range("P1310").Value = 6.53
$ws.Range("Q1310").Value = 30.22
$ws.Range("R1310").Value = 4.72
$ws.Range("S1310").Value = 34
$ws.Range("T1310").Value = 6
$ws.Range("U1310").Value = 33
$ws.Range("V1310").Value = 13

# Row 1311
$ws.Range("A1311").Value = "N3 J16 Bourgoin 14/02/2031"
$ws.Range("B1311").Value = 46067
$ws.Range("C1311").Value = "Global"
$ws.Range("D1311").Value = "M"
$ws.Range("E1311").Value = "Nathanael Beta"
$ws.Range("F1311").Value = "left forward"
$ws.Range("G1311").Value = "00:14:43"
$ws.Range("H1311").Value = 1.76
$ws.Range("I1311").Value = 0.44
$ws.Range("J1311").Value = 1.31
$ws.Range("K1311").Value = 0.23
$ws.Range("L1311").Value = 0.11
$ws.Range("M1311").Value = 0.08
$ws.Range("N1311").Value = 0.03
$ws.Range("O1311").Value = 5
$ws.Range("P1311").Value = 7.11
$ws.Range("Q1311").Value = 32.11
$ws.Range("R1311").Value = 5.09
$ws.Range("S1311").Value = 8
$ws.Range("T1311").Value = 4
$ws.Range("U1311").Value = 6
$ws.Range("V1311").Value = 1

# Row 1312
$ws.Range("A1312").Value = "N3 J16 Bourgoin 14/02/2032"
$ws.Range("B1312").Value = 46067
$ws.Range("C1312").Value = "Global"
$ws.Range("D1312").Value = "M"
$ws.Range("E1312").Value = "Naim Ighbane"
$ws.Range("F1312").Value = "center back"
$ws.Range("G1312").Value = "01:38:10"
$ws.Range("H1312").Value = 10.27
$ws.Range("I1312").Value = 1.53
$ws.Range("J1312").Value = 8.72
$ws.Range("K1312").Value = 1.21
$ws.Range("L1312").Value = 0.28
$ws.Range("M1312").Value = 0.06
$ws.Range("N1312").Value = 0
$ws.Range("O1312").Value = 6
$ws.Range("P1312").Value = 6.17
$ws.Range("Q1312").Value = 28.31
$ws.Range("R1312").Value = 4.31
$ws.Range("S1312").Value = 33
$ws.Range("T1312").Value = 2
$ws.Range("U1312").Value = 23
$ws.Range("V1312").Value = 10

# Row 1313
$ws.Range("A1313").Value = "N3 J16 Bourgoin 14/02/2033"
$ws.Range("B1313").Value = 46067
$ws.Range("C1313").Value = "Global"
$ws.Range("D1313").Value = "M"
$ws.Range("E1313").Value = "Naim Dhib"
$ws.Range("F1313").Value = "center midfield"
$ws.Range("G1313").Value = "01:32:43"
$ws.Range("H1313").Value = 9.74
$ws.Range("I1313").Value = 1.85
$ws.Range("J1313").Value = 7.87
$ws.Range("K1313").Value = 1.25
$ws.Range("L1313").Value = 0.47
$ws.Range("M1313").Value = 0.16
$ws.Range("N1313").Value = 0
$ws.Range("O1313").Value = 8
$ws.Range("P1313").Value = 6.27
$ws.Range("Q1313").Value = 29.33
$ws.Range("R1313").Value = 4.43
$ws.Range("S1313").Value = 46
$ws.Range("T1313").Value = 8
$ws.Range("U1313").Value = 49
$ws.Range("V1313").Value = 14

# Row 1314
$ws.Range("A1314").Value = "N3 J16 Bourgoin 14/02/2034"
$ws.Range("B1314").Value = 46067
$ws.Range("C1314").Value = "Global"
$ws.Range("D1314").Value = "M"
$ws.Range("E1314").Value = "Yoann Martelat"
$ws.Range("F1314").Value = "center midfield"
$ws.Range("G1314").Value = "01:39:05"
$ws.Range("H1314").Value = 12.7
$ws.Range("I1314").Value = 3.28
$ws.Range("J1314").Value = 9.39
$ws.Range("K1314").Value = 2.54
$ws.Range("L1314").Value = 0.71
$ws.Range("M1314").Value = 0.06
$ws.Range("N1314").Value = 0
$ws.Range("O1314").Value = 5
$ws.Range("P1314").Value = 7.61
$ws.Range("Q1314").Value = 27.87
$ws.Range("R1314").Value = 4.46
$ws.Range("S1314").Value = 44
$ws.Range("T1314").Value = 1
$ws.Range("U1314").Value = 49
$ws.Range("V1314").Value = 15

# Row 1315
$ws.Range("A1315").Value = "N3 J16 Bourgoin 14/02/2035"
$ws.Range("B1315").Value = 46067
$ws.Range("C1315").Value = "Global"
$ws.Range("D1315").Value = "M"
$ws.Range("E1315").Value = "Maé Clavel"
$ws.Range("F1315").Value = "left back"
$ws.Range("G1315").Value = "01:39:52"
$ws.Range("H1315").Value = 11.75
$ws.Range("I1315").Value = 1.95
$ws.Range("J1315").Value = 9.77
$ws.Range("K1315").Value = 1.4
$ws.Range("L1315").Value = 0.49
$ws.Range("M1315").Value = 0.08
$ws.Range("N1315").Value = 0
$ws.Range("O1315").Value = 6
$ws.Range("P1315").Value = 7
$ws.Range("Q1315").Value = 27.42
$ws.Range("R1315").Value = 4.69
$ws.Range("S1315").Value = 49
$ws.Range("T1315").Value = 5
$ws.Range("U1315").Value = 41
$ws.Range("V1315").Value = 10

# Row 1316
$ws.Range("A1316").Value = "N3 J16 Bourgoin 14/02/2036"
$ws.Range("B1316").Value = 46067
$ws.Range("C1316").Value = "Global"
$ws.Range("D1316").Value = "M"
$ws.Range("E1316").Value = "Theo Owono"
$ws.Range("F1316").Value = "center midfield"
$ws.Range("G1316").Value = "01:08:46"
$ws.Range("H1316").Value = 8.03
$ws.Range("I1316").Value = 1.72
$ws.Range("J1316").Value = 6.29
$ws.Range("K1316").Value = 1.24
$ws.Range("L1316").Value = 0.37
$ws.Range("M1316").Value = 0.14
$ws.Range("N1316").Value = 0
$ws.Range("O1316").Value = 6
$ws.Range("P1316").Value = 6.95
$ws.Range("Q1316").Value = 29.06
$ws.Range("R1316").Value = 4.4
$ws.Range("S1316").Value = 43
$ws.Range("T1316").Value = 9
$ws.Range("U1316").Value = 39
$ws.Range("V1316").Value = 14

# Row 1317
$ws.Range("A1317").Value = "N3 J16 Bourgoin 14/02/2037"
$ws.Range("B1317").Value = 46067
$ws.Range("C1317").Value = "Global"
$ws.Range("D1317").Value = "M"
$ws.Range("E1317").Value = "Mehdi Boussaid"
$ws.Range("F1317").Value = "right foward"
$ws.Range("G1317").Value = "01:24:23"
$ws.Range("H1317").Value = 9.53
$ws.Range("I1317").Value = 2.04
$ws.Range("J1317").Value = 7.46
$ws.Range("K1317").Value = 1.11
$ws.Range("L1317").Value = 0.7
$ws.Range("M1317").Value = 0.26
$ws.Range("N1317").Value = 0
$ws.Range("O1317").Value = 20
$ws.Range("P1317").Value = 6.66
$ws.Range("Q1317").Value = 29.88
$ws.Range("R1317").Value = 4.32
$ws.Range("S1317").Value = 37
$ws.Range("T1317").Value = 7
$ws.Range("U1317").Value = 23
$ws.Range("V1317").Value = 14

# --- Selection + scroll position, matching the final sheetView in the diff ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1280
$win.ScrollColumn = 1
$ws.Range("A1306:A1317").Select()
